$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 430 (same market/product
# as the existing row 430, but with a newer date), pushing the previous
# rows 430-472 down to 431-473.
$ws.Rows("430:430").Insert(-4121)  # xlShiftDown

$ws.Range("A430").Value = 5
$ws.Range("B430").Value = "Macroferia Regional de Talca"
$ws.Range("C430").Value = "Maule"
$ws.Range("D430").Value = 45132
$ws.Range("E430").Value = 7
$ws.Range("F430").Value = 100112009
$ws.Range("G430").Value = "Acelga"
$ws.Range("H430").Value = "Sin especificar"
$ws.Range("I430").Value = "Primera"
$ws.Range("J430").Value = 500
$ws.Range("K430").Value = 1800
$ws.Range("L430").Value = 1800
$ws.Range("M430").Value = 1800
$ws.Range("N430").Value = "`$/docena de atados (4 kilos)"
$ws.Range("O430").Value = "Región del Maule"
$ws.Range("P430").Value = 450
$ws.Range("Q430").Value = 4
$ws.Range("R430").Value = "Hortaliza"
